# edit.ps1 - applies the "gh-pages output generated at 456a3b4" update to
# 杭州-漫展信息.xlsx
#
# Summary of changes:
#  1. Sheet "展览" (index 1): update "想去人数" (column F) counts for many rows.
#  2. Sheet "演出" (index 2): insert a new row (row 26) for the new event
#     "杭州·世界经典原版音乐剧《猫》CATS", pushing the two existing rows down
#     by one (old row26 -> row27, old row27 -> row28). The running index in
#     column A keeps following the row-1 pattern used throughout the sheet.
#  3. Sheet "本地生活" (index 3): update column F counts for two rows.
#  4. Sheet "全部类型" (index 4): update column F counts for the rows that
#     mirror the other sheets (this sheet is a separate flattened/aggregated
#     list, so the row numbers differ from the source sheets).
#
# NOTE: this runtime's PowerShell-style interpreter does not propagate
# mutations made to COM objects that were passed into a user-defined
# function via `param(...)` - the edits silently vanish. To stay safe,
# every Excel object access below is written out inline rather than
# wrapped in helper functions.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet "展览" - update column F values
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)

$wsExpo.Cells.Item(2, 6).Value = 762
$wsExpo.Cells.Item(3, 6).Value = 14255
$wsExpo.Cells.Item(4, 6).Value = 14352
$wsExpo.Cells.Item(7, 6).Value = 5897
$wsExpo.Cells.Item(9, 6).Value = 575
$wsExpo.Cells.Item(13, 6).Value = 1552
$wsExpo.Cells.Item(14, 6).Value = 441
$wsExpo.Cells.Item(16, 6).Value = 1206
$wsExpo.Cells.Item(17, 6).Value = 1839
$wsExpo.Cells.Item(20, 6).Value = 2293
$wsExpo.Cells.Item(21, 6).Value = 567
$wsExpo.Cells.Item(22, 6).Value = 813
$wsExpo.Cells.Item(23, 6).Value = 3333
$wsExpo.Cells.Item(25, 6).Value = 313
$wsExpo.Cells.Item(26, 6).Value = 2407
$wsExpo.Cells.Item(27, 6).Value = 597
$wsExpo.Cells.Item(28, 6).Value = 118
$wsExpo.Cells.Item(30, 6).Value = 1795
$wsExpo.Cells.Item(31, 6).Value = 1071
$wsExpo.Cells.Item(32, 6).Value = 1399
$wsExpo.Cells.Item(33, 6).Value = 103
$wsExpo.Cells.Item(34, 6).Value = 150
$wsExpo.Cells.Item(35, 6).Value = 4861
$wsExpo.Cells.Item(36, 6).Value = 4864
$wsExpo.Cells.Item(38, 6).Value = 159
$wsExpo.Cells.Item(39, 6).Value = 676
$wsExpo.Cells.Item(41, 6).Value = 3301
$wsExpo.Cells.Item(42, 6).Value = 44
$wsExpo.Cells.Item(43, 6).Value = 923
$wsExpo.Cells.Item(44, 6).Value = 341
$wsExpo.Cells.Item(45, 6).Value = 109
$wsExpo.Cells.Item(46, 6).Value = 85
$wsExpo.Cells.Item(47, 6).Value = 4431
$wsExpo.Cells.Item(48, 6).Value = 592
$wsExpo.Cells.Item(49, 6).Value = 295

# ---------------------------------------------------------------------
# 2. Sheet "演出" - insert the new "猫 CATS" row at row 26
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)

# Shift rows 26:27 down to 27:28 (xlShiftDown = -4121,
# xlFormatFromLeftOrAbove = 0, so the new row inherits row 26's old
# formatting).
$wsShow.Rows.Item(26).Insert(-4121, 0)

# The Insert() call can leave the brand-new row 26 with a freshly minted
# style index for column A instead of reusing the existing one; re-sync
# it by copying the formatting from row 27 (the row that used to be row
# 26, so it still carries the original style index).
$wsShow.Cells.Item(27, 1).Copy() | Out-Null
$wsShow.Cells.Item(26, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Column A is a simple running index equal to (row number - 1); restore
# that invariant for all three affected rows.
$wsShow.Cells.Item(26, 1).Value = 25
$wsShow.Cells.Item(27, 1).Value = 26
$wsShow.Cells.Item(28, 1).Value = 27

# Fill in the new row's data (force text interpretation for the
# date-like strings so they are not auto-converted into date serials,
# then restore the default "Normal" style so no stray number-format
# index is left on the cell).
$wsShow.Cells.Item(26, 2).NumberFormat = "@"
$wsShow.Cells.Item(26, 2).Value = "2024-12-10"
$wsShow.Cells.Item(26, 2).Style = "Normal"

$wsShow.Cells.Item(26, 3).NumberFormat = "@"
$wsShow.Cells.Item(26, 3).Value = "杭州·世界经典原版音乐剧《猫》CATS"
$wsShow.Cells.Item(26, 3).Style = "Normal"

$wsShow.Cells.Item(26, 4).NumberFormat = "@"
$wsShow.Cells.Item(26, 4).Value = "杭州市江干区新业路39号 杭州大剧院"
$wsShow.Cells.Item(26, 4).Style = "Normal"

$wsShow.Cells.Item(26, 5).NumberFormat = "@"
$wsShow.Cells.Item(26, 5).Value = "2024.12.10 19:30-12.15 21:50"
$wsShow.Cells.Item(26, 5).Style = "Normal"

$wsShow.Cells.Item(26, 6).Value = 0
$wsShow.Cells.Item(26, 7).Value = 880

$wsShow.Cells.Item(26, 8).NumberFormat = "@"
$wsShow.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89401"
$wsShow.Cells.Item(26, 8).Style = "Normal"

$wsShow.Cells.Item(26, 9).NumberFormat = "@"
$wsShow.Cells.Item(26, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/aOThG0qq1721123117451.jpeg"
$wsShow.Cells.Item(26, 9).Style = "Normal"

# ---------------------------------------------------------------------
# 3. Sheet "本地生活" - update column F values
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Cells.Item(2, 6).Value = 7589
$wsLocal.Cells.Item(4, 6).Value = 797

# ---------------------------------------------------------------------
# 4. Sheet "全部类型" - update column F values
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)

$wsAll.Cells.Item(2, 6).Value = 762
$wsAll.Cells.Item(4, 6).Value = 797
$wsAll.Cells.Item(6, 6).Value = 14255
$wsAll.Cells.Item(9, 6).Value = 5897
$wsAll.Cells.Item(14, 6).Value = 1552
$wsAll.Cells.Item(15, 6).Value = 441
$wsAll.Cells.Item(16, 6).Value = 1206
$wsAll.Cells.Item(17, 6).Value = 1839
$wsAll.Cells.Item(20, 6).Value = 567
$wsAll.Cells.Item(21, 6).Value = 3333
$wsAll.Cells.Item(22, 6).Value = 313
$wsAll.Cells.Item(23, 6).Value = 597
$wsAll.Cells.Item(25, 6).Value = 1795
$wsAll.Cells.Item(28, 6).Value = 1399
$wsAll.Cells.Item(30, 6).Value = 103
$wsAll.Cells.Item(31, 6).Value = 150
$wsAll.Cells.Item(32, 6).Value = 4861
$wsAll.Cells.Item(33, 6).Value = 4864
$wsAll.Cells.Item(36, 6).Value = 159
$wsAll.Cells.Item(37, 6).Value = 676
$wsAll.Cells.Item(39, 6).Value = 3301
$wsAll.Cells.Item(40, 6).Value = 44
$wsAll.Cells.Item(41, 6).Value = 923
$wsAll.Cells.Item(42, 6).Value = 341
$wsAll.Cells.Item(44, 6).Value = 85
$wsAll.Cells.Item(45, 6).Value = 4431
$wsAll.Cells.Item(46, 6).Value = 592
$wsAll.Cells.Item(47, 6).Value = 295

$wb.Save()
